# Abbreviate the "grupo" (diagnostic group) labels used throughout the
# departamento/grupo table, shorten them to compact codes, widen column A
# to fit the (unchanged) department names, and leave the selection on the
# grupo column that was just edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace each full diagnostic-group label with its abbreviated form.
# These labels repeat once per department (Arauca, Casanare, Meta,
# Vichada), so a sheet-wide Replace updates every occurrence in one shot.
$ws.Cells.Replace("Accidentes de transporte", "Accid. Transp")
$ws.Cells.Replace("Agresiones", "Agres.")
$ws.Cells.Replace("Consumo de sustancias psicoactivas", "Consum.Sust.Psicoact.")
$ws.Cells.Replace("Esquizofrenia Trastornos esquizotípicos y delirantes", "Esquizofrenia")
$ws.Cells.Replace("Lesiones autoinfligidas", "Lesiones Autoinf.")
$ws.Cells.Replace("Retraso mental", "Retr. Mental")
$ws.Cells.Replace("Síndromes del comportamiento  asociados a alteraciones fisiológicas y factores físicos", "Síndr. Alterac. Fisiológ.Fact.Físicos")
$ws.Cells.Replace("Trastornos (afectivos) del estado de ánimo", "Trast. Afect Estad.Animo")
$ws.Cells.Replace("Trastornos de la personalidad y comportamiento en adultos", "Trast.Person.Comp.Adultos")
$ws.Cells.Replace("Trastornos del desarrollo psicológico", "Trast. Desarrollo Psico.")
$ws.Cells.Replace("Trastornos habituales en la niñez y en la adolescencia", "Trast.Habit. Niñez-Adolesc")
$ws.Cells.Replace("Trastornos mentales orgánicos, incluidos los sintomáticos,", "Trast.Ment.Orgán. Sintomát.")
$ws.Cells.Replace("Trastornos neuróticos, trastornos relacionados con el estrés y somatomorfos", "Trast. Neurót. Estrés y Somatom.")

# Column A ("departamento") now needs more room than the default, so size
# it explicitly to fit the longest label.
$ws.Columns.Item(1).ColumnWidth = 24.109375

# Leave the selection on the grupo column (B2:B53) that was just edited,
# matching where the user's cursor ended up.
$ws.Range("B2:B53").Select()
